# Daily attendance processing - reorder "Recorded By" (column G) entries
# Rule: for each comma-separated list of recorders, move the last entry to
# the front of the list (most-recent-first), keeping the remaining entries
# in their original relative order. Single-entry cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    if ($text -notmatch ",") {
        continue
    }

    $parts = $text -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $count = $trimmed.Count
    if ($count -le 1) {
        continue
    }

    $lastItem = $trimmed[$count - 1]
    $rest = $trimmed[0..($count - 2)]

    $newParts = @($lastItem) + $rest
    $newText = [string]::Join(", ", $newParts)

    $cell.Value = $newText
}
